# Update Leave Card 12/22/2023 10:59 AM
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("CONVERTION")

# CONVERTION sheet: enter the number of days (15) used by the
# "late enrollment"/conversion lookup table (J3), which drives
# J4/K3/L3 recalculation.
$ws2.Range("J3").Value = 15

# Sheet1 table row 80 (period row): the period marker moves from
# 8/1/2023 to 8/15/2023, and the EARNED value for the closed period
# (8/1-8/14) is recorded as 0.667.
$ws1.Range("A80").Value = 45153
$ws1.Range("C80").Value = 0.667

# Clear out the future/unused period date placeholders in rows 81-131
# (they were auto-filled and are no longer needed now that the table
# has advanced).
$ws1.Range("A81:A131").ClearContents()
